# ConsolidatedSPO06-23.xlsx — "First Successful Run Through"
#
# 1. Fix two mis-entered ratings on "Staff Performance Overview" (Danielle
#    Mai / Matthew Young had "NA" typed into the Rat column; replace with
#    a text "0" so the column keeps behaving like text, matching the rest
#    of the column).
# 2. Remove the four stray rows (20-23) that had leaked onto the overview
#    sheet — that raw per-employee detail now lives on its own sheet.
# 3. Add a new "Consolidated Data" sheet after the overview sheet and
#    populate it (R2:AT5) with the consolidated metrics for those four
#    employees.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Rat column fixes (E8, E11): "NA" -> "0" (kept as text) -----------
$ws1.Range("E8").NumberFormat = "@"
$ws1.Range("E8").Value = "0"

$ws1.Range("E11").NumberFormat = "@"
$ws1.Range("E11").Value = "0"

# --- 2. Drop the stray rows 20-23 (Chrissy Cummings, Danielle Mai, --------
#        Jasmine Saiz, Karen Trevizo) that were left on the overview sheet
$ws1.Rows("20:23").Delete()

# --- 3. Add the "Consolidated Data" sheet right after the overview sheet -
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Consolidated Data"

# Populate R2:AT5 with the consolidated per-employee detail
$ws2.Range("R2").Value = "Chrissy Cummings"
$ws2.Range("S2").Value = 3
$ws2.Range("T2").Value = 3
$ws2.Range("U2").Value = 1
$ws2.Range("V2").Value = 4.2
$ws2.Range("W2").Value = 251
$ws2.Range("X2").Value = 270.29
$ws2.Range("Y2").Value = 0
$ws2.Range("Z2").Value = 0
$ws2.Range("AA2").Value = 0
$ws2.Range("AB2").Value = 0
$ws2.Range("AC2").Value = 251
$ws2.Range("AD2").Value = 270.29
$ws2.Range("AE2").Value = 83.67
$ws2.Range("AF2").Value = 90.09999999999999
$ws2.Range("AG2").Value = 26
$ws2.Range("AH2").Value = 21
$ws2.Range("AI2").Value = 5
$ws2.Range("AJ2").Value = 4.8
$ws2.Range("AK2").Value = 1843
$ws2.Range("AL2").Value = 1985.86
$ws2.Range("AM2").Value = 0
$ws2.Range("AN2").Value = 0
$ws2.Range("AO2").Value = 0
$ws2.Range("AP2").Value = 0
$ws2.Range("AQ2").Value = 1843
$ws2.Range("AR2").Value = 1985.86
$ws2.Range("AS2").Value = 70.88
$ws2.Range("AT2").Value = 76.38

$ws2.Range("R3").Value = "Danielle Mai"
$ws2.Range("S3").Value = 8
$ws2.Range("T3").Value = 7
$ws2.Range("U3").Value = 1
$ws2.Range("V3").Value = 0
$ws2.Range("W3").Value = 571
$ws2.Range("X3").Value = 614.88
$ws2.Range("Y3").Value = 0
$ws2.Range("Z3").Value = 0
$ws2.Range("AA3").Value = 20
$ws2.Range("AB3").Value = 21.53
$ws2.Range("AC3").Value = 591
$ws2.Range("AD3").Value = 636.41
$ws2.Range("AE3").Value = 73.88
$ws2.Range("AF3").Value = 79.55
$ws2.Range("AG3").Value = 1
$ws2.Range("AH3").Value = 0
$ws2.Range("AI3").Value = 0
$ws2.Range("AJ3").Value = 0
$ws2.Range("AK3").Value = 0
$ws2.Range("AL3").Value = 0
$ws2.Range("AM3").Value = 0
$ws2.Range("AN3").Value = 0
$ws2.Range("AO3").Value = 10
$ws2.Range("AP3").Value = 10.78
$ws2.Range("AQ3").Value = 10
$ws2.Range("AR3").Value = 10.78
$ws2.Range("AS3").Value = 10
$ws2.Range("AT3").Value = 10.78

$ws2.Range("R4").Value = "Jasmine Saiz"
$ws2.Range("S4").Value = 5
$ws2.Range("T4").Value = 5
$ws2.Range("U4").Value = 0
$ws2.Range("V4").Value = 4.2
$ws2.Range("W4").Value = 410
$ws2.Range("X4").Value = 441.51
$ws2.Range("Y4").Value = 0
$ws2.Range("Z4").Value = 0
$ws2.Range("AA4").Value = 0
$ws2.Range("AB4").Value = 0
$ws2.Range("AC4").Value = 410
$ws2.Range("AD4").Value = 441.51
$ws2.Range("AE4").Value = 82
$ws2.Range("AF4").Value = 88.3
$ws2.Range("AG4").Value = 51
$ws2.Range("AH4").Value = 22
$ws2.Range("AI4").Value = 10
$ws2.Range("AJ4").Value = 4.9
$ws2.Range("AK4").Value = 3697
$ws2.Range("AL4").Value = 3983.6
$ws2.Range("AM4").Value = 0
$ws2.Range("AN4").Value = 0
$ws2.Range("AO4").Value = 0
$ws2.Range("AP4").Value = 0
$ws2.Range("AQ4").Value = 3697
$ws2.Range("AR4").Value = 3983.6
$ws2.Range("AS4").Value = 72.48999999999999
$ws2.Range("AT4").Value = 78.11

$ws2.Range("R5").Value = "Karen Trevizo"
$ws2.Range("S5").Value = 19
$ws2.Range("T5").Value = 18
$ws2.Range("U5").Value = 1
$ws2.Range("V5").Value = 4.6
$ws2.Range("W5").Value = 1376
$ws2.Range("X5").Value = 1481.76
$ws2.Range("Y5").Value = 0
$ws2.Range("Z5").Value = 0
$ws2.Range("AA5").Value = 14
$ws2.Range("AB5").Value = 15.08
$ws2.Range("AC5").Value = 1390
$ws2.Range("AD5").Value = 1496.84
$ws2.Range("AE5").Value = 73.16
$ws2.Range("AF5").Value = 78.78
$ws2.Range("AG5").Value = 8
$ws2.Range("AH5").Value = 3
$ws2.Range("AI5").Value = 0
$ws2.Range("AJ5").Value = 4.8
$ws2.Range("AK5").Value = 580
$ws2.Range("AL5").Value = 624.96
$ws2.Range("AM5").Value = 0
$ws2.Range("AN5").Value = 0
$ws2.Range("AO5").Value = 0
$ws2.Range("AP5").Value = 0
$ws2.Range("AQ5").Value = 580
$ws2.Range("AR5").Value = 624.96
$ws2.Range("AS5").Value = 72.5
$ws2.Range("AT5").Value = 78.12
